$d = $word.ActiveDocument

# --- 1) Remove the four paragraphs dropped entirely by the edit ---
#     (indices are in the original, pre-edit document; deleting
#      top-down keeps each index valid as the ones below vanish)
$d.Paragraphs.Item(11).Range.Delete() | Out-Null
$d.Paragraphs.Item(10).Range.Delete() | Out-Null
$d.Paragraphs.Item(9).Range.Delete() | Out-Null
$d.Paragraphs.Item(8).Range.Delete() | Out-Null

# --- 2) Rewrite the text of every paragraph that survives ---
$d.Content.Find.Execute("המאמר היומי של מייק - 14.03.25", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר היומי של מייק - 13.03.25", 2) | Out-Null
$d.Content.Find.Execute("A Survey on Kolmogorov-Arnold Network", $true, $false, $false, $false, $false, $true, 1, $false, "SLIM: Let LLM Learn More and Forget Less with Soft LoRA and Identity Mixture", 2) | Out-Null
$d.Content.Find.Execute("מבוא:", $true, $false, $false, $false, $false, $true, 1, $false, "האם לפעמים קורה לכם שאתם מתחילים לקרוא את המאמר וככל שאתם מתקדמים ומתעמקים בו הוא מתחיל להיראות פחות ופחות טוב. לי זה לפעמים קורה עם אוכל אבל שם יותר קל לי להפסיק לאכול מאשר לקרוא מאמר. אז יאללה, אסקור אותו קצרות אך אל תצפו רבות…", 2) | Out-Null
$d.Content.Find.Execute(" זוכרים את KANs? שזה קיצור של Kolmogorov-Arnold Networks שעשה הרבה רעש בזמנו אך הבאז הלך ודעך עם הזמן. מתברר שיצאו לא מעט מחקרים בנושא המרתק הזה. המאמר דן בהרחבות ושינויים שונים לארכיטקטורת ה-KAN הבסיסית. אלה כוללים התאמות לניתוח סדרות עתיות, לעיבוד דאטה גרפי ולפתרון משוואות דיפרנציאליות. שינויים אלה כוללים לרוב שילוב של רכיבים מיוחדים או אילוצים בתוך ה-KAN במטרה להתמודד טוב יותר עם הדרישות הספציפיות של דומיינים אלה.", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר מציע שיטה להלביש ערבוב של מומחים או MoE על LoRa. נזכיר ש-LoRa היא שיטת פיין טיון של רשתות נוירונים שבהם אנו לא מאמנים את כל משקולות המודל אלא רק מטריצות תוספות בעלת ראנק נמוך. MoE היא שיטה להורדה של העומס החישובי בטרנספורמרים כאשר אנו מחלקים את המטריצות בשכבת FFN של הטרנספורמרים לתת-מטריצות (מומחים) כאשר כל פעם לטוקן נתון אנו מפעילים רק חלק מהמומחים. שכבת ניתוב (routing layer) מחשבות את הציון של כל מומחים ובדרך כלל אנו בוחרים k מומחים בעלי ציון הגבוה ביותר (top-k).", 2) | Out-Null
$d.Content.Find.Execute("רשתות קולמוגורוב-ארנולד מייצגות שינוי פרדיגמה בתכנון רשתות נוירונים, המבוססות על מעבר מפונקציות אקטיבציה קבועות לקראת פונקציות הניתנות ללמידה הנקראות b-splines. הדבר שאב השראה ממשפט הייצוג של קולמוגורוב-ארנולד, הטוען שכל פונקציה רציפה של משתנים מרובים ניתנת לייצוג כהרכבה של פונקציות של משתנה אחד. באמצעות שימוש בפונקציות המיוצגות על ידי ספליינים(שילוב של פולינומים באינטרוול ספוי), KANs מציעות גמישות משופרת ופוטנציאל לדיוק גבוה יותר בקירוב פונקציות. דבר מוביל ל-interpretability משופר של המודל, מכיוון שניתן לנתח ביותר קלות את הפונקציות החד-משתניות שנלמדו.", $true, $false, $false, $false, $false, $true, 1, $false, "אז המחברים משדכים LoRA עם MoE וזה בדיוק מה שמשך את עיניי. המאמר מציע להחליף LoRA רגיל עם כמה מומחי של LoRA שחלקם הינם מטריצות מראנק 0 או פשוט מטריצות אפסים. לטענת המאמר לא תמיד צריך להפעיל את LoRa. מומחי ה-LoRa נבחרים על ידי רשת ניתוב בדומה ל-MoE הסטנדרתי. עבור כל טוקן נבחרים K מומחים (בינם גם מומחי זהות) בעלי ציונים הגבוהים ביותר. שימו לב שבמאמר יש כמה שגיאות בנוסחאות המחשבים את התוצאה של המנגנון המוצע. ", 2) | Out-Null
$d.Content.Find.Execute("רשתות KANs לדומיינים שונים:", $true, $false, $false, $false, $false, $true, 1, $false, "לאחר מכן המאמר מציע שיטה לשכלול הציונים של שכבת הניתוב בהתבסס על הסטטיסטיקות של הדאטהסט עליו בוצע הפיינטיון עם השיטה. סטטיסטיקה במקרה הזה מחושבת על המצבים החבויים של הרשת המחושבים על הדאטה של הפיין טיון (אופן החישוב המדויק לא מוגדר בצורה ברורה ולדעתי יש שגיאות בנוסחאות המגדירות אותו). המחברים מציעים לקלסטר את המצבים החבויים האלו לקלסטרים שמספרם כנראה שווה למספר הטוקנים בפרומפט (מוגדר כקבוע במאמר ועבור סדרות קצרות יותר משתמשים בטוקני ה-padding). ", 2) | Out-Null
$d.Content.Find.Execute("כעת נתאר כמה הרחבות של KAN לדומיינים שונים. לניתוח סדרות עתיות, רשתות KAN זמניות (T-KANs) משלבות מנגנוני זיכרון, בדומה ל-RNNs ו-LSTM, לטיפול יעיל בסדרות אלו ובתלויות לטווח ארוך שבהן, ומדגימות ביצועים מעולים במשימות חיזוי רב-שלבי(multi-step forecasting). בנוסף, שינויים כמו מנגנונים חיבורים gated, בדומה LSTM ו-GRU, מאפשרים ל-KANs להתאים באופן דינמי פונקציות אקטיבציה (ספליין בגדול* בהתבסס על מורכבות המשימה, משפרים יעילות מבלי לדרוש רגולריזציה נרחבת.", $true, $false, $false, $false, $false, $true, 1, $false, "מרכזי הקלסטרים מתעדכנים במהלך הפיין טיון (כל פרומפט הקלט משויך לקלסטר הקרוב ביותר ואז מרכז הקלסטר מחושב מחדש). במהלך האינפרנס פרומפט הקלט משויך לקלסטר הקרוב ביותר (מרחק ריבוע) ואז ציוני המומחים המופקים על ידי שכבת הניתוב עבור מומחי הזהות מוזזים במקדם שעולה אם המרחק לקלסטר הקרוב עולה כאשר הציונים למומחי LoRA האחרים נותרים ללא שינוי. נציין שמרכזי הקלסטרים לא מתעדכנים במהלך האינפרנס.", 2) | Out-Null
$d.Content.Find.Execute("בדאטה הגרפי, KANs מבוססות גרף (GKANs) פותחו לשיפור סיווג צמתים semi-supervised על ידי שיפור זרימת מידע בין צמתים, עולות בביצועיהן על רשתות קונבולוציה גרפיות מסורתיות (GCNs). ארכיטקטורות מבוססות KAN אלה משפרות את למידת ייצוג הצמתים ומשפרות את דיוק מודלי הרגרסיה בגרפים העולות ברשתות חברתיות וכימיה מולקולרית. GCNs פועלות על ידי צבירה ושינוי חוזרים של מידע תכונות משכונות מקומיות בתוך גרף, ותופסות ביעילות הן תכונות צמתים והן טופולוגיית גרף. עם זאת, GCNs מסתמכות על פילטרי קונבולוציה קבועים, המגבילים את הגמישות שלהן בטיפול בגרפים מורכבים והטרוגניים. כדי להתמודד עם מגבלה זו, GKAN מציג שתי ארכיטקטורות עיקריות: ארכיטקטורה 1, המצרפת תכונות צמתים לפני יישום שכבות KAN, מאפשרת לפונקציות אקטיבציה הניתנות ללמידה לתפוס יחסים מקומיים מורכבים, וארכיטקטורה 2, הממקמת שכבות KAN בין הטמעות צמתים בכל שכבה לפני הצבירה, מאפשרת התאמה דינמית לשינויים במבנה הגרף. שיפור זה מאפשר ל-GKANs להסתגל באופן דינמי לשינויים במבנה הגרף, ומספק גישה יותר אדפטיבית ללמידה מבוססת גרף.", $true, $false, $false, $false, $false, $true, 1, $false, "לבסוף המאמר מציע דרך לשלב כמה MoE עם LoRa עבור כמה משימות פיין טיון שונות אבל אחרי שגיליתי טעיות גם בפרק הזה, ויתרתי….", 2) | Out-Null
$d.Content.Find.Execute("https://arxiv.org/abs/2411.06078", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/pdf/2410.07739", 2) | Out-Null
